# Fix suspendu mail content issue
# Append a new "YAYA TATA " row (row 11) to the "Etat Virement" table,
# mirroring the existing rows' structure/types.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text / identifier columns -> force text so number-like strings
# (CIN/IF, account numbers, contract numbers, ...) are not reinterpreted
# as numeric values.
$ws.Range("A11").Value = "YAYA TATA "

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "KL365695"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "354654613156464166848965"

# D11 (Agence) is blank for this row, same as D7/D8 for other "Supervision" rows.
$ws.Range("D11").Value = ""

$ws.Range("E11").Value = "BP"
$ws.Range("F11").Value = "Supervision"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "569/SUP 9999"

$ws.Range("H11").Value = "mensuelle"

# Numeric amount columns.
$ws.Range("I11").Value = 80000
$ws.Range("J11").Value = 12000
$ws.Range("K11").Value = 68000
